$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 3972.6365
$ws.Range("I51").Value = 2999.75
$ws.Range("J51").Value = 4528.5713
$ws.Range("K51").Value = 2999.75
$ws.Range("L51").Value = 4528.5713
$ws.Range("M51").Value = -2515.75
$ws.Range("N51").Value = -5496.5713
# Row 107
$ws.Range("H107").Value = 546.06665
$ws.Range("I107").Value = 368.57144
$ws.Range("J107").Value = 701.375
$ws.Range("K107").Value = 368.57144
$ws.Range("L107").Value = 701.375
$ws.Range("M107").Value = 1551.42856
$ws.Range("N107").Value = -4541.375
# Row 129
$ws.Range("H129").Value = 1086.8628
$ws.Range("J129").Value = 1244.0714
$ws.Range("L129").Value = 3732.2142
$ws.Range("N129").Value = -13732.2142
# Row 137
$ws.Range("H137").Value = 1315.2433
$ws.Range("I137").Value = 1308.7
$ws.Range("J137").Value = 1343.2858
$ws.Range("K137").Value = 3926.1
$ws.Range("L137").Value = 4029.8574
$ws.Range("M137").Value = -1376.1
$ws.Range("N137").Value = -9129.857400000001
# Row 138
$ws.Range("H138").Value = 2086.4788
$ws.Range("I138").Value = 1144.279
$ws.Range("J138").Value = 3533.4285
$ws.Range("K138").Value = 3432.837
$ws.Range("L138").Value = 10600.2855
$ws.Range("M138").Value = 1707.163
$ws.Range("N138").Value = -20880.2855

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11459.028
$ws.Range("I32").Value = 11575.6455
$ws.Range("J32").Value = 10555.25
$ws.Range("K32").Value = 11575.6455
$ws.Range("L32").Value = 10555.25
$ws.Range("M32").Value = -11288.6455
$ws.Range("N32").Value = -11129.25
# Row 44
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976
$ws.Range("M44").ClearContents()
# Row 109
$ws.Range("H109").Value = 38624.875
$ws.Range("J109").Value = 38624.875
$ws.Range("L109").Value = 38624.875
$ws.Range("N109").Value = -41398.875
# Row 110
$ws.Range("H110").Value = 1301.25
$ws.Range("I110").Value = 999
$ws.Range("J110").Value = 1402
$ws.Range("K110").Value = 999
$ws.Range("L110").Value = 1402
$ws.Range("M110").Value = 1046
$ws.Range("N110").Value = -5492
# Row 132
$ws.Range("H132").Value = 4390.977
$ws.Range("I132").Value = 5924.24
$ws.Range("J132").Value = 2373.5264
$ws.Range("K132").Value = 17772.72
$ws.Range("L132").Value = 7120.5792
$ws.Range("M132").Value = -15242.72
$ws.Range("N132").Value = -12180.5792

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 7939607
$ws.Range("I105").Value = 8406437
$ws.Range("K105").Value = 8406437
$ws.Range("M105").Value = -8404690
# Row 108
$ws.Range("H108").Value = 25122.625
$ws.Range("J108").Value = 25122.625
$ws.Range("L108").Value = 25122.625
$ws.Range("N108").Value = -32802.625

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 9999.546
$ws.Range("J4").Value = 9999.546
$ws.Range("L4").Value = 9999.546
$ws.Range("N4").Value = -10223.546
# Row 16
$ws.Range("H16").Value = 1426.1
$ws.Range("I16").Value = 1410
$ws.Range("J16").Value = 1450.25
$ws.Range("K16").Value = 1410
$ws.Range("L16").Value = 1450.25
$ws.Range("M16").Value = -1123
$ws.Range("N16").Value = -2024.25
# Row 31
$ws.Range("H31").Value = 2492.5518
$ws.Range("I31").Value = 1649.2727
$ws.Range("J31").Value = 5142.857
$ws.Range("K31").Value = 1649.2727
$ws.Range("L31").Value = 5142.857
$ws.Range("M31").Value = -1354.2727
$ws.Range("N31").Value = -5732.857
# Row 34
$ws.Range("H34").Value = 2492.5518
$ws.Range("I34").Value = 1649.2727
$ws.Range("J34").Value = 5142.857
$ws.Range("K34").Value = 1649.2727
$ws.Range("L34").Value = 5142.857
$ws.Range("M34").Value = -1447.2727
$ws.Range("N34").Value = -5546.857
# Row 94
$ws.Range("H94").Value = 2015.5333
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 2312.0908
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 2312.0908
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -3214.0908
# Row 113
$ws.Range("H113").Value = 1426.1
$ws.Range("I113").Value = 1410
$ws.Range("J113").Value = 1450.25
$ws.Range("K113").Value = 1410
$ws.Range("L113").Value = 1450.25
$ws.Range("M113").Value = 760
$ws.Range("N113").Value = -5790.25
# Row 132
$ws.Range("H132").Value = 339574.25
$ws.Range("I132").Value = 398828.6
$ws.Range("K132").Value = 1196485.8
$ws.Range("M132").Value = -1193955.8

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 30785.715
$ws.Range("J93").Value = 30785.715
$ws.Range("L93").Value = 30785.715
$ws.Range("N93").Value = -34529.715
# Row 102
$ws.Range("H102").Value = 3827.182
$ws.Range("I102").Value = 4499.75
$ws.Range("K102").Value = 4499.75
$ws.Range("M102").Value = -2877.75
# Row 109
$ws.Range("H109").Value = 9480.75
$ws.Range("J109").Value = 9480.75
$ws.Range("L109").Value = 9480.75
$ws.Range("N109").Value = -11560.75
# Row 126
$ws.Range("H126").Value = 3822.2856
$ws.Range("I126").Value = 3845.7778
$ws.Range("J126").Value = 3780
$ws.Range("K126").Value = 11537.3334
$ws.Range("L126").Value = 11340
$ws.Range("M126").Value = -9067.3334
$ws.Range("N126").Value = -16280
# Row 132
$ws.Range("H132").Value = 2015.3414
$ws.Range("I132").Value = 1357
$ws.Range("J132").Value = 3044
$ws.Range("K132").Value = 4071
$ws.Range("L132").Value = 9132
$ws.Range("M132").Value = -1541
$ws.Range("N132").Value = -14192

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 7877.5557
$ws.Range("I132").Value = 8983.333000000001
$ws.Range("K132").Value = 26949.999
$ws.Range("M132").Value = -24419.999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 500750
$ws.Range("I5").Value = 500750
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500750
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -500638
$ws.Range("N5").ClearContents()
# Row 132
$ws.Range("H132").Value = 2084.8823
$ws.Range("I132").Value = 1359.65
$ws.Range("J132").Value = 3120.9285
$ws.Range("K132").Value = 4078.95
$ws.Range("L132").Value = 9362.7855
$ws.Range("M132").Value = -1548.95
$ws.Range("N132").Value = -14422.7855
